$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Loans")

$newDate = "Wed May 23 16:27:54 GMT-03:00 2018"
$newExpiration = "Thu May 24 16:27:54 GMT-03:00 2018"

for ($r = 5; $r -le 7; $r++) {
    $ws.Cells.Item($r, 4).Value = $newDate
    $ws.Cells.Item($r, 5).Value = $newExpiration
}
